# Comercializadora del Agro de Limarí - Frutilla
# A new weekly price-report date block (44951) is inserted before the
# existing row 451 block, pushing the rest of the table down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 451:453 (existing rows 451.. shift down to 454..)
$ws.Rows("451:453").Insert()

# Shared/common column values for every data row in this sheet
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100101
$producto   = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "$/bandeja 7 kilos"
$origen      = "Provincia de Melipilla"
$kgUnidad    = 7
$fecha       = 44951

$rows = @(
    @{ Row = 451; Calidad = "Especial"; Volumen = 560; PMin = 13000; PMax = 14000; PProm = 13500; PKg = 1929 },
    @{ Row = 452; Calidad = "Primera";  Volumen = 600; PMin = 11000; PMax = 12000; PProm = 11500; PKg = 1643 },
    @{ Row = 453; Calidad = "Segunda";  Volumen = 500; PMin = 9000;  PMax = 10000; PProm = 9500;  PKg = 1357 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
